$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "reason" header in C1, matching the style used by the other headers (A1/B1)
$ws.Range("C1").Value = "reason"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Row 2 (id=3, score=100) - add reason text
$ws.Range("C2").Value = "The job as an NLP engineer requires skills in NLP, PyTorch, Computer Vision, and Python, which align with your expertise in computer vision and NLP. Additionally, the job involves analyzing and preprocessing large-scale text data, which is a direct match to your experience in object detection and classification. The high score indicates a strong suitability for your skillset."

# Row 3 (id=1) - score 66 -> 64, add reason text
$ws.Range("B3").Value = 64
$ws.Range("C3").Value = "The job of SDE Intern with a score of 64 is moderately suitable for the candidate. The candidate has experience in ReactJS, JavaScript, and web development, which aligns with the required skills. However, their projects do not directly involve MongoDB and NodeJS, which may have contributed to the moderate score. Overall, the candidate's experience makes them a potential fit for the role, but they may need some additional training in MongoDB and NodeJS."

# Row 4 (id=2) - score 65 -> 63, add reason text
$ws.Range("B4").Value = 63
$ws.Range("C4").Value = "The job as a Frontend Engineer Intern requires skills in ReactJS, JavaScript, CSS, Frontend Development, and NextJS, which align with the candidate's experience in projects like the 'Website for the Literature Society' and 'LLMGuard'. The moderate score of 63 suggests that while the candidate possesses relevant skills, there may be other candidates with stronger qualifications or experiences."
